# TP3/MaDUM.xlsx — add two more coverage matrices next to the existing one:
#   * I10:P15  -> exact duplicate of the A10:H15 table (same header/rows)
#   * A17:G22  -> the A10:H15 table again, but without the "unzip_tree" column
# plus a green highlight on the "build_codebook" (J) column of the new block,
# and on L11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- duplicate the whole existing table (rows 10-15, cols A-H) into I:P ----
$ws.Range("A10:H15").Copy($ws.Range("I10:P15"))

# ---- build the reduced table (rows 17-22) missing the unzip_tree column ----
# columns A-E (index, __init__, build_codebook, from_string, encode_tree) stay put
$ws.Range("A10:E15").Copy($ws.Range("A17:E22"))
# columns G-H (huffman_encode, huffman_decode) slide left into F-G, dropping F (unzip_tree)
$ws.Range("G10:H15").Copy($ws.Range("F17:G22"))

# ---- highlight the build_codebook column (green) on the new I:P block ----
$green = 5296274  # RGB(146, 208, 80) == 0xFF92D050
$ws.Range("J11").Interior.Color = $green
$ws.Range("J12").Interior.Color = $green
$ws.Range("J13").Interior.Color = $green
$ws.Range("J14").Interior.Color = $green
$ws.Range("J15").Interior.Color = $green
$ws.Range("L11").Interior.Color = $green

# ---- widen the data columns to cover the new columns up to P ----
$ws.Range("B1:P1").ColumnWidth = 14.6

# ---- view state: selection + zoom ----
$ws.Range("K12").Select() | Out-Null
$excel.ActiveWindow.Zoom = 85

Write-Output "edit applied"
